$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "com.alm.wiley.controller.Customer"
$ws.Range("B1").Value = "com.alm.wiley.controller.drools"

[void]$ws.Range("B2").Select()
